$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update comment strings: "lam" -> "egp"
$ws.Range("L5").Value = "Début de l’effet du confinement (egp 1.12 → 1)"
$ws.Range("L11").Value = "Fin supposée de l’épidémie (egp=1)"

# Update predicted growth-rate values in column B (rows 4-10)
$ws.Range("B4").Value = 1.05
$ws.Range("B5").Value = 1.03
$ws.Range("B6").Value = 1.03
$ws.Range("B7").Value = 1.02
$ws.Range("B8").Value = 1.02
$ws.Range("B9").Value = 1.01
$ws.Range("B10").Value = 1.005

# Update active cell selection
$ws.Range("L12").Select()
